$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 386.66666
$ws.Range("I12").Value = 330
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 330
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -160
$ws.Range("N12").Value = -840
$ws.Range("H28").Value = 4027.2173
$ws.Range("I28").Value = 1052.5385
$ws.Range("K28").Value = 1052.5385
$ws.Range("M28").Value = -567.5385000000001
$ws.Range("H64").Value = 4564.846
$ws.Range("J64").Value = 4344.375
$ws.Range("L64").Value = 4344.375
$ws.Range("N64").Value = -4840.375
$ws.Range("H67").Value = 4564.846
$ws.Range("J67").Value = 4344.375
$ws.Range("L67").Value = 4344.375
$ws.Range("N67").Value = -6060.375
$ws.Range("H80").Value = 5604423.5
$ws.Range("I80").Value = 3269873.5
$ws.Range("J80").Value = 7938974
$ws.Range("K80").Value = 9809620.5
$ws.Range("L80").Value = 23816922
$ws.Range("M80").Value = -9808622.5
$ws.Range("N80").Value = -23818918
$ws.Range("H83").Value = 5604423.5
$ws.Range("I83").Value = 3269873.5
$ws.Range("J83").Value = 7938974
$ws.Range("K83").Value = 29428861.5
$ws.Range("L83").Value = 71450766
$ws.Range("M83").Value = -29423869.5
$ws.Range("N83").Value = -71460750
$ws.Range("H86").Value = 4878
$ws.Range("I86").Value = 2320.75
$ws.Range("K86").Value = 2320.75
$ws.Range("M86").Value = -1197.75
$ws.Range("H88").Value = 5126.7144
$ws.Range("I88").Value = 3715.3333
$ws.Range("J88").Value = 6185.25
$ws.Range("K88").Value = 3715.3333
$ws.Range("L88").Value = 6185.25
$ws.Range("M88").Value = -3309.3333
$ws.Range("N88").Value = -6997.25
$ws.Range("H89").Value = 4878
$ws.Range("I89").Value = 2320.75
$ws.Range("K89").Value = 11603.75
$ws.Range("M89").Value = -5987.75
$ws.Range("H91").Value = 5126.7144
$ws.Range("I91").Value = 3715.3333
$ws.Range("J91").Value = 6185.25
$ws.Range("K91").Value = 3715.3333
$ws.Range("L91").Value = 6185.25
$ws.Range("M91").Value = -2311.3333
$ws.Range("N91").Value = -8993.25
$ws.Range("H100").Value = 4621.95
$ws.Range("I100").Value = 2924.5557
$ws.Range("K100").Value = 2924.5557
$ws.Range("M100").Value = -2383.5557
$ws.Range("H116").Value = 17228.7
$ws.Range("I116").Value = 6900
$ws.Range("J116").Value = 18376.334
$ws.Range("K116").Value = 6900
$ws.Range("L116").Value = 18376.334
$ws.Range("M116").Value = -3458
$ws.Range("N116").Value = -25260.334
$ws.Range("H137").Value = 20837642
$ws.Range("I137").Value = 45457296
$ws.Range("K137").Value = 136371888
$ws.Range("M137").Value = -136369338
$ws.Range("H138").Value = 10402.5625
$ws.Range("J138").Value = 13902.2
$ws.Range("L138").Value = 41706.60000000001
$ws.Range("N138").Value = -51986.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1950.1666
$ws.Range("I32").Value = 1776.4048
$ws.Range("K32").Value = 1776.4048
$ws.Range("M32").Value = -1489.4048
$ws.Range("H43").Value = 37591.57
$ws.Range("J43").Value = 36950.2
$ws.Range("L43").Value = 36950.2
$ws.Range("N43").Value = -37576.2
$ws.Range("H48").Value = 335161.34
$ws.Range("J48").Value = 335161.34
$ws.Range("L48").Value = 335161.34
$ws.Range("N48").Value = -335929.34
$ws.Range("H61").Value = 78889920
$ws.Range("I61").Value = 140000800
$ws.Range("K61").Value = 140000800
$ws.Range("M61").Value = -140000588
$ws.Range("H102").Value = 2052.15
$ws.Range("I102").Value = 1864
$ws.Range("J102").Value = 2491.1667
$ws.Range("K102").Value = 1864
$ws.Range("L102").Value = 2491.1667
$ws.Range("M102").Value = -242
$ws.Range("N102").Value = -5735.1667
$ws.Range("H132").Value = 4171358
$ws.Range("I132").Value = 4913.3184
$ws.Range("K132").Value = 14739.9552
$ws.Range("M132").Value = -12209.9552
$ws.Range("H136").Value = 78889920
$ws.Range("I136").Value = 140000800
$ws.Range("K136").Value = 420002400
$ws.Range("M136").Value = -419999850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3484.4
$ws.Range("I20").Value = 3105.5
$ws.Range("K20").Value = 3105.5
$ws.Range("M20").Value = -2858.5
$ws.Range("H86").Value = 22157.412
$ws.Range("I86").Value = 40821.53
$ws.Range("K86").Value = 40821.53
$ws.Range("M86").Value = -39698.53
$ws.Range("H89").Value = 22157.412
$ws.Range("I89").Value = 40821.53
$ws.Range("K89").Value = 204107.65
$ws.Range("M89").Value = -198491.65
$ws.Range("H105").Value = 463625.06
$ws.Range("I105").Value = 759170.9399999999
$ws.Range("J105").Value = 6872.364
$ws.Range("K105").Value = 759170.9399999999
$ws.Range("L105").Value = 6872.364
$ws.Range("M105").Value = -757423.9399999999
$ws.Range("N105").Value = -10366.364
$ws.Range("H134").Value = 6668668
$ws.Range("I134").Value = 1835.4166
$ws.Range("K134").Value = 5506.2498
$ws.Range("M134").Value = -2971.2498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 14288907
$ws.Range("I16").Value = 20003270
$ws.Range("K16").Value = 20003270
$ws.Range("M16").Value = -20002983
$ws.Range("H86").Value = 16320.333
$ws.Range("I86").Value = 15860.375
$ws.Range("K86").Value = 15860.375
$ws.Range("M86").Value = -14737.375
$ws.Range("H89").Value = 16320.333
$ws.Range("I89").Value = 15860.375
$ws.Range("K89").Value = 79301.875
$ws.Range("M89").Value = -73685.875
$ws.Range("H107").Value = 1536.3334
$ws.Range("I107").Value = 1343.9706
$ws.Range("J107").Value = 2130.9092
$ws.Range("K107").Value = 1343.9706
$ws.Range("L107").Value = 2130.9092
$ws.Range("M107").Value = 576.0293999999999
$ws.Range("N107").Value = -5970.9092
$ws.Range("H113").Value = 14288907
$ws.Range("I113").Value = 20003270
$ws.Range("K113").Value = 20003270
$ws.Range("M113").Value = -20001100
$ws.Range("H132").Value = 2818.4285
$ws.Range("I132").Value = 2755.7896
$ws.Range("K132").Value = 8267.3688
$ws.Range("M132").Value = -5737.3688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 75002824
$ws.Range("I80").Value = 2994.5
$ws.Range("J80").Value = 100002770
$ws.Range("K80").Value = 2994.5
$ws.Range("L80").Value = 100002770
$ws.Range("M80").Value = -1996.5
$ws.Range("N80").Value = -100004766
$ws.Range("H83").Value = 75002824
$ws.Range("I83").Value = 2994.5
$ws.Range("J83").Value = 100002770
$ws.Range("K83").Value = 14972.5
$ws.Range("L83").Value = 500013850
$ws.Range("M83").Value = -9980.5
$ws.Range("N83").Value = -500023834
$ws.Range("H132").Value = 11937933
$ws.Range("I132").Value = 5330.6665
$ws.Range("K132").Value = 15991.9995
$ws.Range("M132").Value = -13461.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 81000.5
$ws.Range("I53").Value = 81000.5
$ws.Range("K53").Value = 81000.5
$ws.Range("M53").Value = -80482.5
$ws.Range("H68").Value = 5720835.5
$ws.Range("I68").Value = 7480231.5
$ws.Range("J68").Value = 2798
$ws.Range("K68").Value = 7480231.5
$ws.Range("L68").Value = 2798
$ws.Range("M68").Value = -7479482.5
$ws.Range("N68").Value = -4296
$ws.Range("H71").Value = 5720835.5
$ws.Range("I71").Value = 7480231.5
$ws.Range("J71").Value = 2798
$ws.Range("K71").Value = 37401157.5
$ws.Range("L71").Value = 13990
$ws.Range("M71").Value = -37397413.5
$ws.Range("N71").Value = -21478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1444.7858
$ws.Range("I81").Value = 1171.3077
$ws.Range("K81").Value = 2342.6154
$ws.Range("M81").Value = -1281.6154
$ws.Range("H84").Value = 1444.7858
$ws.Range("I84").Value = 1171.3077
$ws.Range("K84").Value = 11713.077
$ws.Range("M84").Value = -6409.077000000001
$ws.Range("H122").Value = 3579.1
$ws.Range("I122").Value = 3058
$ws.Range("K122").Value = 9174
$ws.Range("M122").Value = -6724
